$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Test Cases")

# --- Row 33: E33 result flips from SKIP to PASS ---
$ws.Range("E33").Value = "PASS"

# --- Row 49: E49 result flips from SKIP to PASS ---
$ws.Range("E49").Value = "PASS"

# --- New row 50: FollowUnfollowPostsAuthor / OPQA-427 ---
$ws.Range("A50").Value = "FollowUnfollowPostsAuthor"
$ws.Range("C50").Value = "Veirfy that the user is able to follow the author of the post directly from the post"
$ws.Hyperlinks.Add($ws.Range("B50"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-427", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-427")
$ws.Range("B50").Value = "OPQA-427"
$ws.Range("A49:E49").Copy()
$ws.Range("A50:E50").PasteSpecial(-4122)
$ws.Range("D50").Value = "Y"
$ws.Range("E50").Value = "PASS"

# --- New row 51: CommentOnUsersOwnPost / OPQA-377 ---
$ws.Range("A51").Value = "CommentOnUsersOwnPost"
$ws.Range("C51").Value = "Verify that the user is able to comment on the post a user authored themselves."
$ws.Range("B51").Value = "OPQA-377"
$ws.Range("A49:E49").Copy()
$ws.Range("A51:E51").PasteSpecial(-4122)
$ws.Range("D51").Value = "Y"
$ws.Range("E51").Value = "PASS"

# --- Update selection to the new last row ---
[void]$ws.Range("A51:E51").Select()
